$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.220.07'
$ws.Range('E2').Value = '  -1.75%  '

$ws.Range('D3').Value = '2.182.27'
$ws.Range('E3').Value = '  -1.86%  '

$ws.Range('E4').Value = '  -0.11%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.83'
$ws.Range('E5').Value = '  -0.86%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.613'
$ws.Range('E6').Value = '  -2.96%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '66.52'
$ws.Range('E7').Value = '  -6.64%  '

$ws.Range('E8').Value = '  -0.02%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.585'
$ws.Range('E9').Value = '  -2.45%  '

$ws.Range('B10').Value = 'OKB'
$ws.Range('C10').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '59.29'
$ws.Range('E10').Value = '  +1.72%  '

$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '37.10'
$ws.Range('E11').Value = '  -8.95%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0934'
$ws.Range('E12').Value = '  -2.87%  '

$ws.Range('E13').Value = '  -0.74%  '

$ws.Range('E14').Value = '  -5.48%  '

$ws.Range('D15').Value = '2.508.09'
$ws.Range('E15').Value = '  -1.81%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.30'
$ws.Range('E16').Value = '  -4.77%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.846'
$ws.Range('E17').Value = '  -3.77%  '

$ws.Range('D18').Value = '2.172.57'
$ws.Range('E18').Value = '  -2.44%  '

$ws.Range('D19').Value = '41.138.75'
$ws.Range('E19').Value = '  -1.69%  '

$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').Value = '  -2.21%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.61'
$ws.Range('E21').Value = '  -1.61%  '

$ws.Range('E22').Value = '  -3.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '230.18'
$ws.Range('E23').Value = '  -2.34%  '

$ws.Range('E24').Value = '  -3.85%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.83'
$ws.Range('E25').Value = '  -6.20%  '

$ws.Range('E26').Value = '  -0.02%  '

$ws.Range('E27').Value = '  -2.90%  '

$ws.Range('E28').Value = '  -5.64%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.67'
$ws.Range('E29').Value = '  -3.51%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '166.46'
$ws.Range('E30').Value = '  -2.44%  '

$ws.Range('E31').Value = '  -7.67%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.14'
$ws.Range('E32').Value = '  -3.48%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.119'
$ws.Range('E33').Value = '  -2.62%  '

$ws.Range('E34').Value = '  +2.95%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0753'
$ws.Range('E35').Value = '  +1.75%  '

$ws.Range('E36').Value = '  -2.34%  '

$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.09'
$ws.Range('E37').Value = '  +1.28%  '

$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.52'
$ws.Range('E38').Value = '  -4.28%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '25.01'
$ws.Range('E39').Value = '  -4.84%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0304'
$ws.Range('E40').Value = '  -0.94%  '

$ws.Range('B41').Value = 'LidoDAOToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.21'
$ws.Range('E41').Value = '  -3.27%  '

$ws.Range('B42').Value = 'FTXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.41'
$ws.Range('E42').Value = '  +12.39%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.55'
$ws.Range('E43').Value = '  -6.57%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '11.37'
$ws.Range('E44').Value = '  -10.79%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '60.12'
$ws.Range('E45').Value = '  -7.67%  '

$ws.Range('E46').Value = '  -2.81%  '

$ws.Range('E47').Value = '  -8.07%  '

$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0990'
$ws.Range('E48').Value = '  -3.57%  '

$ws.Range('B49').Value = 'BinanceUSD'
$ws.Range('C49').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  -0.25%  '

$ws.Range('E50').Value = '  -3.42%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.30'
$ws.Range('E51').Value = '  -7.39%  '
